# Update gh-pages to output generated at 456a3b4
# Applies the diff between the previous scrape and the new scrape of
# 广州-漫展信息.xlsx across the 4 worksheets: 展览, 演出, 本地生活, 全部类型.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    # Force the value to be stored as literal text even when it looks like
    # a date (e.g. "2024-10-05"), then strip the residual quote-prefix style
    # so the cell's style index matches a plain/unstyled cell.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 41
$ws1.Range("F4").Value = 25866
$ws1.Range("F5").Value = 569
$ws1.Range("F6").Value = 241
$ws1.Range("F7").Value = 577
$ws1.Range("F8").Value = 167
$ws1.Range("F11").Value = 343
$ws1.Range("F12").Value = 203
$ws1.Range("F13").Value = 173
$ws1.Range("F15").Value = 276
Set-TextCell $ws1 "C16" "广州·wio流金序曲乙女同人ONLY展"
$ws1.Range("F16").Value = 337
$ws1.Range("F18").Value = 1477
$ws1.Range("F19").Value = 159

# Insert a new row before row 20 (shifts old rows 20-22 down to 21-23)
$ws1.Rows.Item(20).Insert()

$ws1.Range("A20").Value = 19
$ws1.Range("A20").Font.Bold = $true
$ws1.Range("A20").HorizontalAlignment = -4108
$ws1.Range("A20").VerticalAlignment = -4160
$ws1.Range("A20").Borders.LineStyle = 1

Set-TextCell $ws1 "B20" "2024-10-05"
Set-TextCell $ws1 "C20" "广州·南部动漫节"
Set-TextCell $ws1 "D20" "东沙大道16号 广州健康方舟"
Set-TextCell $ws1 "E20" "2024.10.05 10:00-10.06 17:00"
$ws1.Range("F20").Value = 4
$ws1.Range("G20").Value = 60
Set-TextCell $ws1 "H20" "https://show.bilibili.com/platform/detail.html?id=90923"
Set-TextCell $ws1 "I20" "//i0.hdslb.com/bfs/openplatform/202408/RMKuGZYX1722580310264.jpeg"

# Old row 20 (now row 21), old row 21 (now row 22), old row 22 (now row 23)
$ws1.Range("F21").Value = 411
$ws1.Range("F22").Value = 90
$ws1.Range("F23").Value = 113

# ---------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

Set-TextCell $ws2 "G2" "已售罄"
$ws2.Range("F6").Value = 78
$ws2.Range("F7").Value = 33
$ws2.Range("F8").Value = 101
$ws2.Range("F9").Value = 101
$ws2.Range("F10").Value = 422
$ws2.Range("F18").Value = 18

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 4904
$ws3.Range("F3").Value = 178

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 41
$ws4.Range("F4").Value = 4904
$ws4.Range("F5").Value = 178
$ws4.Range("F6").Value = 25866
$ws4.Range("F7").Value = 569
Set-TextCell $ws4 "G8" "已售罄"
$ws4.Range("F9").Value = 241
$ws4.Range("F11").Value = 577
$ws4.Range("F14").Value = 167
$ws4.Range("F15").Value = 78
$ws4.Range("F16").Value = 78
$ws4.Range("F17").Value = 33
$ws4.Range("F18").Value = 101
$ws4.Range("F19").Value = 101
$ws4.Range("F20").Value = 422
$ws4.Range("F24").Value = 343
$ws4.Range("F25").Value = 203
$ws4.Range("F26").Value = 173
$ws4.Range("F29").Value = 276
Set-TextCell $ws4 "C32" "广州·wio流金序曲乙女同人ONLY展"
$ws4.Range("F32").Value = 337
$ws4.Range("F35").Value = 1477
$ws4.Range("F36").Value = 159

# Insert a new row before row 38 (shifts old rows 38-44 down to 39-45)
$ws4.Rows.Item(38).Insert()

$ws4.Range("A38").Value = 37
$ws4.Range("A38").Font.Bold = $true
$ws4.Range("A38").HorizontalAlignment = -4108
$ws4.Range("A38").VerticalAlignment = -4160
$ws4.Range("A38").Borders.LineStyle = 1

Set-TextCell $ws4 "B38" "2024-10-05"
Set-TextCell $ws4 "C38" "广州·南部动漫节"
Set-TextCell $ws4 "D38" "东沙大道16号 广州健康方舟"
Set-TextCell $ws4 "E38" "2024.10.05 10:00-10.06 17:00"
$ws4.Range("F38").Value = 4
$ws4.Range("G38").Value = 60
Set-TextCell $ws4 "H38" "https://show.bilibili.com/platform/detail.html?id=90923"
Set-TextCell $ws4 "I38" "//i0.hdslb.com/bfs/openplatform/202408/RMKuGZYX1722580310264.jpeg"

$ws4.Range("F39").Value = 411
$ws4.Range("F40").Value = 90
$ws4.Range("F41").Value = 113
